$d = $word.ActiveDocument

$d.Content.Find.Execute("11×62=", $true, $false, $false, $false, $false, $true, 1, $false, "67×28=", 2) | Out-Null
$d.Content.Find.Execute("41×28=", $true, $false, $false, $false, $false, $true, 1, $false, "31×70=", 2) | Out-Null
$d.Content.Find.Execute("51×46=", $true, $false, $false, $false, $false, $true, 1, $false, "30×56=", 2) | Out-Null
$d.Content.Find.Execute("28×13=", $true, $false, $false, $false, $false, $true, 1, $false, "73×25=", 2) | Out-Null
$d.Content.Find.Execute("25×21=", $true, $false, $false, $false, $false, $true, 1, $false, "42×20=", 2) | Out-Null
$d.Content.Find.Execute("59×74=", $true, $false, $false, $false, $false, $true, 1, $false, "78×96=", 2) | Out-Null
$d.Content.Find.Execute("47×34=", $true, $false, $false, $false, $false, $true, 1, $false, "70×14=", 2) | Out-Null
$d.Content.Find.Execute("19×97=", $true, $false, $false, $false, $false, $true, 1, $false, "84×58=", 2) | Out-Null
$d.Content.Find.Execute("60×93=", $true, $false, $false, $false, $false, $true, 1, $false, "92×39=", 2) | Out-Null
$d.Content.Find.Execute("99×83=", $true, $false, $false, $false, $false, $true, 1, $false, "69×54=", 2) | Out-Null
$d.Content.Find.Execute("50×19=", $true, $false, $false, $false, $false, $true, 1, $false, "11×78=", 2) | Out-Null
$d.Content.Find.Execute("43×22=", $true, $false, $false, $false, $false, $true, 1, $false, "98×63=", 2) | Out-Null
$d.Content.Find.Execute("17×37=", $true, $false, $false, $false, $false, $true, 1, $false, "96×54=", 2) | Out-Null
$d.Content.Find.Execute("15×13=", $true, $false, $false, $false, $false, $true, 1, $false, "91×84=", 2) | Out-Null
$d.Content.Find.Execute("60×26=", $true, $false, $false, $false, $false, $true, 1, $false, "96×65=", 2) | Out-Null
$d.Content.Find.Execute("79×94=", $true, $false, $false, $false, $false, $true, 1, $false, "66×40=", 2) | Out-Null
$d.Content.Find.Execute("63×85=", $true, $false, $false, $false, $false, $true, 1, $false, "81×42=", 2) | Out-Null
$d.Content.Find.Execute("93×50=", $true, $false, $false, $false, $false, $true, 1, $false, "48×100=", 2) | Out-Null
$d.Content.Find.Execute("86×21=", $true, $false, $false, $false, $false, $true, 1, $false, "40×50=", 2) | Out-Null
$d.Content.Find.Execute("39×63=", $true, $false, $false, $false, $false, $true, 1, $false, "50×52=", 2) | Out-Null
$d.Content.Find.Execute("36×31=", $true, $false, $false, $false, $false, $true, 1, $false, "31×69=", 2) | Out-Null
$d.Content.Find.Execute("22×15=", $true, $false, $false, $false, $false, $true, 1, $false, "96×49=", 2) | Out-Null
$d.Content.Find.Execute("69×100=", $true, $false, $false, $false, $false, $true, 1, $false, "27×66=", 2) | Out-Null
$d.Content.Find.Execute("21×23=", $true, $false, $false, $false, $false, $true, 1, $false, "18×16=", 2) | Out-Null
$d.Content.Find.Execute("33×90=", $true, $false, $false, $false, $false, $true, 1, $false, "97×100=", 2) | Out-Null
$d.Content.Find.Execute("21×82=", $true, $false, $false, $false, $false, $true, 1, $false, "36×58=", 2) | Out-Null
$d.Content.Find.Execute("38×70=", $true, $false, $false, $false, $false, $true, 1, $false, "44×40=", 2) | Out-Null
$d.Content.Find.Execute("43×78=", $true, $false, $false, $false, $false, $true, 1, $false, "12×59=", 2) | Out-Null
$d.Content.Find.Execute("32×87=", $true, $false, $false, $false, $false, $true, 1, $false, "75×56=", 2) | Out-Null
$d.Content.Find.Execute("55×87=", $true, $false, $false, $false, $false, $true, 1, $false, "32×40=", 2) | Out-Null
$d.Content.Find.Execute("97×78=", $true, $false, $false, $false, $false, $true, 1, $false, "91×95=", 2) | Out-Null
$d.Content.Find.Execute("73×59=", $true, $false, $false, $false, $false, $true, 1, $false, "60×53=", 2) | Out-Null
$d.Content.Find.Execute("74×42=", $true, $false, $false, $false, $false, $true, 1, $false, "10×94=", 2) | Out-Null
$d.Content.Find.Execute("94×24=", $true, $false, $false, $false, $false, $true, 1, $false, "27×19=", 2) | Out-Null
$d.Content.Find.Execute("27×30=", $true, $false, $false, $false, $false, $true, 1, $false, "99×98=", 2) | Out-Null
$d.Content.Find.Execute("72×77=", $true, $false, $false, $false, $false, $true, 1, $false, "58×48=", 2) | Out-Null
$d.Content.Find.Execute("87×81=", $true, $false, $false, $false, $false, $true, 1, $false, "42×56=", 2) | Out-Null
$d.Content.Find.Execute("27×68=", $true, $false, $false, $false, $false, $true, 1, $false, "58×60=", 2) | Out-Null
$d.Content.Find.Execute("95×79=", $true, $false, $false, $false, $false, $true, 1, $false, "18×31=", 2) | Out-Null
$d.Content.Find.Execute("44×18=", $true, $false, $false, $false, $false, $true, 1, $false, "42×79=", 2) | Out-Null
$d.Content.Find.Execute("98×100=", $true, $false, $false, $false, $false, $true, 1, $false, "92×54=", 2) | Out-Null
$d.Content.Find.Execute("53×99=", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=", 2) | Out-Null
$d.Content.Find.Execute("46×77=", $true, $false, $false, $false, $false, $true, 1, $false, "36×97=", 2) | Out-Null
$d.Content.Find.Execute("20×61=", $true, $false, $false, $false, $false, $true, 1, $false, "68×81=", 2) | Out-Null
$d.Content.Find.Execute("58×13=", $true, $false, $false, $false, $false, $true, 1, $false, "23×39=", 2) | Out-Null
$d.Content.Find.Execute("54×82=", $true, $false, $false, $false, $false, $true, 1, $false, "75×50=", 2) | Out-Null
$d.Content.Find.Execute("56×29=", $true, $false, $false, $false, $false, $true, 1, $false, "85×49=", 2) | Out-Null
$d.Content.Find.Execute("60×19=", $true, $false, $false, $false, $false, $true, 1, $false, "92×93=", 2) | Out-Null
$d.Content.Find.Execute("18×54=", $true, $false, $false, $false, $false, $true, 1, $false, "67×14=", 2) | Out-Null
$d.Content.Find.Execute("79×45=", $true, $false, $false, $false, $false, $true, 1, $false, "69×41=", 2) | Out-Null
$d.Content.Find.Execute("94×69=", $true, $false, $false, $false, $false, $true, 1, $false, "24×24=", 2) | Out-Null
$d.Content.Find.Execute("19×64=", $true, $false, $false, $false, $false, $true, 1, $false, "66×71=", 2) | Out-Null
$d.Content.Find.Execute("46×68=", $true, $false, $false, $false, $false, $true, 1, $false, "78×81=", 2) | Out-Null
$d.Content.Find.Execute("30×21=", $true, $false, $false, $false, $false, $true, 1, $false, "35×39=", 2) | Out-Null
$d.Content.Find.Execute("71×61=", $true, $false, $false, $false, $false, $true, 1, $false, "10×51=", 2) | Out-Null
$d.Content.Find.Execute("68×22=", $true, $false, $false, $false, $false, $true, 1, $false, "36×36=", 2) | Out-Null
$d.Content.Find.Execute("77×39=", $true, $false, $false, $false, $false, $true, 1, $false, "92×46=", 2) | Out-Null
$d.Content.Find.Execute("100×43=", $true, $false, $false, $false, $false, $true, 1, $false, "37×15=", 2) | Out-Null
$d.Content.Find.Execute("96×32=", $true, $false, $false, $false, $false, $true, 1, $false, "64×80=", 2) | Out-Null
$d.Content.Find.Execute("74×69=", $true, $false, $false, $false, $false, $true, 1, $false, "93×94=", 2) | Out-Null
$d.Content.Find.Execute("46×63=", $true, $false, $false, $false, $false, $true, 1, $false, "53×100=", 2) | Out-Null
$d.Content.Find.Execute("89×74=", $true, $false, $false, $false, $false, $true, 1, $false, "72×78=", 2) | Out-Null
$d.Content.Find.Execute("67×55=", $true, $false, $false, $false, $false, $true, 1, $false, "54×22=", 2) | Out-Null
$d.Content.Find.Execute("41×48=", $true, $false, $false, $false, $false, $true, 1, $false, "32×25=", 2) | Out-Null
$d.Content.Find.Execute("48×90=", $true, $false, $false, $false, $false, $true, 1, $false, "21×73=", 2) | Out-Null
$d.Content.Find.Execute("93×78=", $true, $false, $false, $false, $false, $true, 1, $false, "84×73=", 2) | Out-Null
$d.Content.Find.Execute("95×91=", $true, $false, $false, $false, $false, $true, 1, $false, "45×85=", 2) | Out-Null
$d.Content.Find.Execute("22×52=", $true, $false, $false, $false, $false, $true, 1, $false, "54×59=", 2) | Out-Null
$d.Content.Find.Execute("64×18=", $true, $false, $false, $false, $false, $true, 1, $false, "34×94=", 2) | Out-Null
$d.Content.Find.Execute("99×80=", $true, $false, $false, $false, $false, $true, 1, $false, "25×57=", 2) | Out-Null
$d.Content.Find.Execute("19×40=", $true, $false, $false, $false, $false, $true, 1, $false, "13×90=", 2) | Out-Null
$d.Content.Find.Execute("65×65=", $true, $false, $false, $false, $false, $true, 1, $false, "24×72=", 2) | Out-Null
$d.Content.Find.Execute("52×40=", $true, $false, $false, $false, $false, $true, 1, $false, "29×34=", 2) | Out-Null
$d.Content.Find.Execute("29×58=", $true, $false, $false, $false, $false, $true, 1, $false, "72×44=", 2) | Out-Null
$d.Content.Find.Execute("75×58=", $true, $false, $false, $false, $false, $true, 1, $false, "71×29=", 2) | Out-Null
$d.Content.Find.Execute("26×48=", $true, $false, $false, $false, $false, $true, 1, $false, "11×24=", 2) | Out-Null
$d.Content.Find.Execute("60×79=", $true, $false, $false, $false, $false, $true, 1, $false, "92×14=", 2) | Out-Null
$d.Content.Find.Execute("65×82=", $true, $false, $false, $false, $false, $true, 1, $false, "20×35=", 2) | Out-Null
$d.Content.Find.Execute("65×76=", $true, $false, $false, $false, $false, $true, 1, $false, "34×12=", 2) | Out-Null
$d.Content.Find.Execute("11×73=", $true, $false, $false, $false, $false, $true, 1, $false, "74×100=", 2) | Out-Null
$d.Content.Find.Execute("87×26=", $true, $false, $false, $false, $false, $true, 1, $false, "96×25=", 2) | Out-Null
$d.Content.Find.Execute("44×54=", $true, $false, $false, $false, $false, $true, 1, $false, "15×99=", 2) | Out-Null
$d.Content.Find.Execute("24×19=", $true, $false, $false, $false, $false, $true, 1, $false, "11×50=", 2) | Out-Null
$d.Content.Find.Execute("31×54=", $true, $false, $false, $false, $false, $true, 1, $false, "53×77=", 2) | Out-Null
$d.Content.Find.Execute("54×95=", $true, $false, $false, $false, $false, $true, 1, $false, "37×59=", 2) | Out-Null
$d.Content.Find.Execute("34×17=", $true, $false, $false, $false, $false, $true, 1, $false, "98×22=", 2) | Out-Null
$d.Content.Find.Execute("61×93=", $true, $false, $false, $false, $false, $true, 1, $false, "82×91=", 2) | Out-Null
$d.Content.Find.Execute("20×16=", $true, $false, $false, $false, $false, $true, 1, $false, "43×20=", 2) | Out-Null
$d.Content.Find.Execute("15×51=", $true, $false, $false, $false, $false, $true, 1, $false, "63×89=", 2) | Out-Null
$d.Content.Find.Execute("17×86=", $true, $false, $false, $false, $false, $true, 1, $false, "32×59=", 2) | Out-Null
$d.Content.Find.Execute("33×14=", $true, $false, $false, $false, $false, $true, 1, $false, "53×25=", 2) | Out-Null
$d.Content.Find.Execute("46×44=", $true, $false, $false, $false, $false, $true, 1, $false, "60×75=", 2) | Out-Null
$d.Content.Find.Execute("93×96=", $true, $false, $false, $false, $false, $true, 1, $false, "19×20=", 2) | Out-Null
$d.Content.Find.Execute("37×64=", $true, $false, $false, $false, $false, $true, 1, $false, "14×76=", 2) | Out-Null
$d.Content.Find.Execute("42×59=", $true, $false, $false, $false, $false, $true, 1, $false, "95×36=", 2) | Out-Null
$d.Content.Find.Execute("62×96=", $true, $false, $false, $false, $false, $true, 1, $false, "57×18=", 2) | Out-Null
$d.Content.Find.Execute("26×37=", $true, $false, $false, $false, $false, $true, 1, $false, "22×25=", 2) | Out-Null
$d.Content.Find.Execute("16×74=", $true, $false, $false, $false, $false, $true, 1, $false, "65×70=", 2) | Out-Null
$d.Content.Find.Execute("55×97=", $true, $false, $false, $false, $false, $true, 1, $false, "19×95=", 2) | Out-Null
$d.Content.Find.Execute("31×35=", $true, $false, $false, $false, $false, $true, 1, $false, "69×38=", 2) | Out-Null
